$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shift the existing model-run rows (7-10) down to (8-11) so a brand-new
# "Constant q" model run can be written into row 7. Only the results block
# (H:S) moves - the A/B price-history columns and the G run-index column
# are left alone.
# ---------------------------------------------------------------------------
$src = $ws.Range("H7:S10")
$dst = $ws.Range("H8:S11")
$src.Copy($dst)

# Range.Copy only carries the cached values in this engine, so the relative
# formulas in the P ("N parameters" running count) and Q ("Other") columns
# need to be re-written by hand, shifted down one row just like a real
# Excel insert would do.
$ws.Range("P10").Formula = "=P9+1"
$ws.Range("P11").Formula = "=P10+1"
$ws.Range("Q8").Formula  = "=I8-SUM(J8:O8)+S8"
$ws.Range("Q9").Formula  = "=I9-SUM(J9:O9)+S9"
$ws.Range("Q10").Formula = "=I10-SUM(J10:O10)+S10"
$ws.Range("Q11").Formula = "=I11-SUM(J11:O11)+S11"

# ---------------------------------------------------------------------------
# New row 7: "Constant q" model run
# ---------------------------------------------------------------------------
$ws.Range("H7").Value    = "Constant q"
$ws.Range("I7").Value2   = 1396.02
$ws.Range("J7").Value2   = 108.339
$ws.Range("K7").Value2   = 0.170124
$ws.Range("L7").Value2   = 621.18
$ws.Range("M7").Value2   = 572.857
$ws.Range("N7").Value2   = 24.05059
$ws.Range("O7").Value2   = 65.8667
$ws.Range("P7").Value2   = 487
$ws.Range("Q7").Formula  = "=I7-SUM(J7:O7)+S7"
$ws.Range("R7").Value2   = 0
$ws.Range("S7").Value2   = 0.267414

# ---------------------------------------------------------------------------
# Clear the stale "run index" numbers out of column G (rows 7-10 and
# 12-15); the formatting/style of those cells is left in place.
# ---------------------------------------------------------------------------
foreach ($r in 7, 8, 9, 10, 12, 13, 14, 15) {
    $ws.Cells.Item($r, 7).ClearContents()
}

# ---------------------------------------------------------------------------
# Update the visible selection to reflect the newly added block.
# ---------------------------------------------------------------------------
$ws.Range("F16:W31").Select()

$wb.Save()
